# Updated symbol list (crypto price/volume refresh) on Tue Jan 31 03:31:13 UTC 2023.
# All Price/Volume(1h) cells are stored as text (t="inlineStr") in the source
# workbook, so numeric-looking values are written with a leading apostrophe
# to force Excel to keep them as text instead of coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.84"
$ws.Range("E2").Value = "'-1.57%"
$ws.Range("D3").Value = "'38.04"
$ws.Range("E3").Value = "'-3.67%"
$ws.Range("D4").Value = "'5.059"
$ws.Range("E4").Value = "'-1.25%"
$ws.Range("D5").Value = "'0.07767"
$ws.Range("E5").Value = "'-4.89%"
$ws.Range("D6").Value = "'4.353"
$ws.Range("E6").Value = "'-0.28%"
$ws.Range("D7").Value = "'1.907"
$ws.Range("E7").Value = "'-3.52%"
$ws.Range("D8").Value = "'8.197"
$ws.Range("E8").Value = "'-1.61%"
$ws.Range("D9").Value = "'0.9225"
$ws.Range("E9").Value = "'-1.68%"
$ws.Range("D10").Value = "'0.1243"
$ws.Range("E10").Value = "'-5.42%"
$ws.Range("D11").Value = "'0.1876"
$ws.Range("E11").Value = "'-4.80%"
$ws.Range("D12").Value = "'0.08763"
$ws.Range("E12").Value = "'-3.20%"
$ws.Range("D13").Value = "'0.03418"
$ws.Range("E13").Value = "'-2.26%"
$ws.Range("D14").Value = "'0.09706"
$ws.Range("E14").Value = "'-0.37%"
$ws.Range("D15").Value = "'0.001373"
$ws.Range("E15").Value = "'-2.93%"
$ws.Range("D16").Value = "'0.005952"
$ws.Range("E16").Value = "'-1.78%"
$ws.Range("D17").Value = "'3.570"
$ws.Range("E17").Value = "'-2.22%"
$ws.Range("D18").Value = "'3.086"
$ws.Range("E18").Value = "'-6.48%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'5.029"
$ws.Range("E20").Value = "'1.26%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1267"
$ws.Range("E21").Value = "'-3.73%"
$ws.Range("D22").Value = "'0.2618"
$ws.Range("E22").Value = "'5.09%"
$ws.Range("E23").Value = "'5,592.83%"
$ws.Range("D24").Value = "'0.04385"
$ws.Range("E24").Value = "'0.14%"
$ws.Range("D25").Value = "'0.001211"
$ws.Range("E25").Value = "'-2.46%"
$ws.Range("D26").Value = "'0.004257"
$ws.Range("E26").Value = "'-10.71%"
$ws.Range("E27").Value = "'-65.32%"
$ws.Range("D39").Value = "'0.02135"
$ws.Range("E39").Value = "'-3.36%"
$ws.Range("D40").Value = "'0.04996"
$ws.Range("E40").Value = "'-3.70%"
$ws.Range("D41").Value = "'0.007952"
$ws.Range("E41").Value = "'2.42%"
$ws.Range("D42").Value = "'0.01004"
$ws.Range("E42").Value = "'-3.10%"
$ws.Range("D43").Value = "'0.1341"
$ws.Range("E43").Value = "'-4.38%"
$ws.Range("E44").Value = "'-1.96%"
$ws.Range("D45").Value = "'0.008785"
$ws.Range("E45").Value = "'-5.37%"
$ws.Range("D46").Value = "'0.00006474"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.00%"
$ws.Range("D48").Value = "'0.003227"
$ws.Range("E48").Value = "'11.82%"
$ws.Range("E49").Value = "'-0.15%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.00%"
